# "changed direction to web app"
# Recolor the "Transfer" highlight fill from green to sky blue, and push the
# three "Rent" entries' timestamps later in the day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of the cells sharing the green "Transfer" highlight fill - set them all
# to the same new sky-blue color together so they end up sharing one fill.
$transferCells = @("D4", "D7", "D9", "D14", "D18", "D21")
foreach ($addr in $transferCells) {
    $rng = $ws.Range($addr)
    $rng.Interior.Color = 15453831        # RGB(135, 206, 235) = #87CEEB
    $rng.Interior.PatternColor = 15453831
}

# Shift the "Rent" row timestamps (date stays the same, time-of-day moves
# from the early morning to the evening)
$ws.Range("A8").Value = 45382.9414271412
$ws.Range("A15").Value = 45412.9414271412
$ws.Range("A22").Value = 45442.9414271412
